$wb = $excel.ActiveWorkbook

# ----- Sheet 1: LP1912 -----
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 1).Value = 'Última actualización: 10:13:53'
$ws.Cells.Item(3, 1).Value = 'Total filas: 141'
$ws.Cells.Item(48, 1).Value = '05:53:46'
$ws.Cells.Item(48, 2).Value = '07:31'
$ws.Cells.Item(48, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(48, 4).Value = 98
$ws.Cells.Item(48, 5).Value = 'LP1912'
$ws.Cells.Item(49, 1).Value = '07:24:45'
$ws.Cells.Item(49, 2).Value = '07:31'
$ws.Cells.Item(49, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(49, 4).Value = 7
$ws.Cells.Item(49, 5).Value = 'LP1912'
$ws.Cells.Item(116, 1).Value = '10:13:53'
$ws.Cells.Item(116, 2).Value = '10:22'
$ws.Cells.Item(116, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(116, 4).Value = 9
$ws.Cells.Item(116, 5).Value = 'LP1912'
$ws.Cells.Item(117, 1).Value = '09:26:30'
$ws.Cells.Item(117, 2).Value = '10:23'
$ws.Cells.Item(117, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(117, 4).Value = 57
$ws.Cells.Item(117, 5).Value = 'LP1912'
$ws.Cells.Item(118, 1).Value = '08:31:53'
$ws.Cells.Item(118, 2).Value = '10:26'
$ws.Cells.Item(118, 3).Value = '215A_EL PATO'
$ws.Cells.Item(118, 4).Value = 115
$ws.Cells.Item(118, 5).Value = 'LP1912'
$ws.Cells.Item(119, 1).Value = '10:13:53'
$ws.Cells.Item(119, 2).Value = '10:32'
$ws.Cells.Item(119, 3).Value = '10_OLMOS'
$ws.Cells.Item(119, 4).Value = 19
$ws.Cells.Item(119, 5).Value = 'LP1912'
$ws.Cells.Item(120, 1).Value = '10:13:53'
$ws.Cells.Item(120, 2).Value = '10:34'
$ws.Cells.Item(120, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(120, 4).Value = 21
$ws.Cells.Item(120, 5).Value = 'LP1912'
$ws.Cells.Item(121, 1).Value = '10:13:53'
$ws.Cells.Item(121, 2).Value = '10:34'
$ws.Cells.Item(121, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(121, 4).Value = 21
$ws.Cells.Item(121, 5).Value = 'LP1912'
$ws.Cells.Item(122, 1).Value = '08:47:51'
$ws.Cells.Item(122, 2).Value = '10:41'
$ws.Cells.Item(122, 3).Value = '17_ROMERO'
$ws.Cells.Item(122, 4).Value = 114
$ws.Cells.Item(122, 5).Value = 'LP1912'
$ws.Cells.Item(123, 1).Value = '08:55:25'
$ws.Cells.Item(123, 2).Value = '10:42'
$ws.Cells.Item(123, 3).Value = '17_ROMERO'
$ws.Cells.Item(123, 4).Value = 107
$ws.Cells.Item(123, 5).Value = 'LP1912'
$ws.Cells.Item(124, 1).Value = '08:47:51'
$ws.Cells.Item(124, 2).Value = '10:43'
$ws.Cells.Item(124, 3).Value = '14_ABASTO'
$ws.Cells.Item(124, 4).Value = 116
$ws.Cells.Item(124, 5).Value = 'LP1912'
$ws.Cells.Item(125, 1).Value = '10:13:53'
$ws.Cells.Item(125, 2).Value = '10:52'
$ws.Cells.Item(125, 3).Value = '15_ABASTO'
$ws.Cells.Item(125, 4).Value = 39
$ws.Cells.Item(125, 5).Value = 'LP1912'
$ws.Cells.Item(126, 1).Value = '10:13:53'
$ws.Cells.Item(126, 2).Value = '10:56'
$ws.Cells.Item(126, 3).Value = '27_EL RETIRO'
$ws.Cells.Item(126, 4).Value = 43
$ws.Cells.Item(126, 5).Value = 'LP1912'
$ws.Cells.Item(127, 1).Value = '09:26:30'
$ws.Cells.Item(127, 2).Value = '10:57'
$ws.Cells.Item(127, 3).Value = '27_EL RETIRO'
$ws.Cells.Item(127, 4).Value = 91
$ws.Cells.Item(127, 5).Value = 'LP1912'
$ws.Cells.Item(128, 1).Value = '10:13:53'
$ws.Cells.Item(128, 2).Value = '11:01'
$ws.Cells.Item(128, 3).Value = '215C_EL PATO'
$ws.Cells.Item(128, 4).Value = 48
$ws.Cells.Item(128, 5).Value = 'LP1912'
$ws.Cells.Item(129, 1).Value = '09:26:30'
$ws.Cells.Item(129, 2).Value = '11:02'
$ws.Cells.Item(129, 3).Value = '215C_EL PATO'
$ws.Cells.Item(129, 4).Value = 96
$ws.Cells.Item(129, 5).Value = 'LP1912'
$ws.Cells.Item(130, 1).Value = '10:13:53'
$ws.Cells.Item(130, 2).Value = '11:03'
$ws.Cells.Item(130, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(130, 4).Value = 50
$ws.Cells.Item(130, 5).Value = 'LP1912'
$ws.Cells.Item(131, 1).Value = '09:26:30'
$ws.Cells.Item(131, 2).Value = '11:06'
$ws.Cells.Item(131, 3).Value = '16_P MOR-167 Y 521'
$ws.Cells.Item(131, 4).Value = 100
$ws.Cells.Item(131, 5).Value = 'LP1912'
$ws.Cells.Item(132, 1).Value = '10:13:53'
$ws.Cells.Item(132, 2).Value = '11:10'
$ws.Cells.Item(132, 3).Value = '10_OLMOS'
$ws.Cells.Item(132, 4).Value = 57
$ws.Cells.Item(132, 5).Value = 'LP1912'
$ws.Cells.Item(133, 1).Value = '10:13:53'
$ws.Cells.Item(133, 2).Value = '11:12'
$ws.Cells.Item(133, 3).Value = '15_ABASTO'
$ws.Cells.Item(133, 4).Value = 59
$ws.Cells.Item(133, 5).Value = 'LP1912'
$ws.Cells.Item(134, 1).Value = '10:13:53'
$ws.Cells.Item(134, 2).Value = '11:12'
$ws.Cells.Item(134, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(134, 4).Value = 59
$ws.Cells.Item(134, 5).Value = 'LP1912'
$ws.Cells.Item(135, 1).Value = '09:26:30'
$ws.Cells.Item(135, 2).Value = '11:19'
$ws.Cells.Item(135, 3).Value = '86_EST CHICA-ESC AGRARIA'
$ws.Cells.Item(135, 4).Value = 113
$ws.Cells.Item(135, 5).Value = 'LP1912'
$ws.Cells.Item(136, 1).Value = '10:13:53'
$ws.Cells.Item(136, 2).Value = '11:20'
$ws.Cells.Item(136, 3).Value = '26_HERNANDEZ'
$ws.Cells.Item(136, 4).Value = 67
$ws.Cells.Item(136, 5).Value = 'LP1912'
$ws.Cells.Item(137, 1).Value = '09:26:30'
$ws.Cells.Item(137, 2).Value = '11:21'
$ws.Cells.Item(137, 3).Value = '26_HERNANDEZ'
$ws.Cells.Item(137, 4).Value = 115
$ws.Cells.Item(137, 5).Value = 'LP1912'
$ws.Cells.Item(138, 1).Value = '10:13:53'
$ws.Cells.Item(138, 2).Value = '11:26'
$ws.Cells.Item(138, 3).Value = '225_C ROCA-H SUR'
$ws.Cells.Item(138, 4).Value = 73
$ws.Cells.Item(138, 5).Value = 'LP1912'
$ws.Cells.Item(139, 1).Value = '10:13:53'
$ws.Cells.Item(139, 2).Value = '11:32'
$ws.Cells.Item(139, 3).Value = '81_EL PELIGRO'
$ws.Cells.Item(139, 4).Value = 79
$ws.Cells.Item(139, 5).Value = 'LP1912'
$ws.Cells.Item(140, 1).Value = '10:13:53'
$ws.Cells.Item(140, 2).Value = '11:38'
$ws.Cells.Item(140, 3).Value = '10_OLMOS'
$ws.Cells.Item(140, 4).Value = 85
$ws.Cells.Item(140, 5).Value = 'LP1912'
$ws.Cells.Item(141, 1).Value = '10:13:53'
$ws.Cells.Item(141, 2).Value = '11:41'
$ws.Cells.Item(141, 3).Value = '17_ROMERO'
$ws.Cells.Item(141, 4).Value = 88
$ws.Cells.Item(141, 5).Value = 'LP1912'
$ws.Cells.Item(142, 1).Value = '10:13:53'
$ws.Cells.Item(142, 2).Value = '11:51'
$ws.Cells.Item(142, 3).Value = '215B_EL PATO'
$ws.Cells.Item(142, 4).Value = 98
$ws.Cells.Item(142, 5).Value = 'LP1912'
$ws.Cells.Item(143, 1).Value = '10:13:53'
$ws.Cells.Item(143, 2).Value = '11:58'
$ws.Cells.Item(143, 3).Value = '225_GOMEZ'
$ws.Cells.Item(143, 4).Value = 105
$ws.Cells.Item(143, 5).Value = 'LP1912'
$ws.Cells.Item(144, 1).Value = '10:13:53'
$ws.Cells.Item(144, 2).Value = '12:02'
$ws.Cells.Item(144, 3).Value = '84_COLONIA URQUIZA-ESC 49'
$ws.Cells.Item(144, 4).Value = 109
$ws.Cells.Item(144, 5).Value = 'LP1912'
$ws.Cells.Item(145, 1).Value = '10:13:53'
$ws.Cells.Item(145, 2).Value = '12:06'
$ws.Cells.Item(145, 3).Value = '16_P MOR-SANTA ANA'
$ws.Cells.Item(145, 4).Value = 113
$ws.Cells.Item(145, 5).Value = 'LP1912'
$ws.Cells.Item(146, 1).Value = '10:13:53'
$ws.Cells.Item(146, 2).Value = '12:06'
$ws.Cells.Item(146, 3).Value = '14_ABASTO'
$ws.Cells.Item(146, 4).Value = 113
$ws.Cells.Item(146, 5).Value = 'LP1912'

# ----- Sheet 2: LP1912-215 -----
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 1).Value = 'Última actualización: 10:13:53'
$ws.Cells.Item(3, 1).Value = 'Total filas: 20'
$ws.Cells.Item(23, 1).Value = '10:13:53'
$ws.Cells.Item(23, 2).Value = '11:01'
$ws.Cells.Item(23, 3).Value = '215C_EL PATO'
$ws.Cells.Item(23, 4).Value = 48
$ws.Cells.Item(23, 5).Value = 'LP1912'
$ws.Cells.Item(24, 1).Value = '09:26:30'
$ws.Cells.Item(24, 2).Value = '11:02'
$ws.Cells.Item(24, 3).Value = '215C_EL PATO'
$ws.Cells.Item(24, 4).Value = 96
$ws.Cells.Item(24, 5).Value = 'LP1912'
$ws.Cells.Item(25, 1).Value = '10:13:53'
$ws.Cells.Item(25, 2).Value = '11:51'
$ws.Cells.Item(25, 3).Value = '215B_EL PATO'
$ws.Cells.Item(25, 4).Value = 98
$ws.Cells.Item(25, 5).Value = 'LP1912'

# ----- Sheet 3: 6203-6173 -----
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 1).Value = 'Última actualización: 10:13:53'
$ws.Cells.Item(3, 1).Value = 'Total filas: 29'
$ws.Cells.Item(32, 1).Value = '10:13:53'
$ws.Cells.Item(32, 2).Value = '11:13'
$ws.Cells.Item(32, 3).Value = '215C_LA PLATA'
$ws.Cells.Item(32, 4).Value = 60
$ws.Cells.Item(32, 5).Value = 'L6203'
$ws.Cells.Item(33, 1).Value = '09:26:30'
$ws.Cells.Item(33, 2).Value = '11:14'
$ws.Cells.Item(33, 3).Value = '215C_LA PLATA'
$ws.Cells.Item(33, 4).Value = 108
$ws.Cells.Item(33, 5).Value = 'L6203'
$ws.Cells.Item(34, 1).Value = '10:13:53'
$ws.Cells.Item(34, 2).Value = '12:04'
$ws.Cells.Item(34, 3).Value = '215A_LA PLATA'
$ws.Cells.Item(34, 4).Value = 111
$ws.Cells.Item(34, 5).Value = 'L6173'
